$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume/coin data per latest scrape.
# NumberFormat "@" + Style "Normal" keeps numeric-looking text (e.g. "1.00")
# stored as text instead of being coerced to a number, while leaving the
# cell formatting identical to its original (default) style.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '55.732.21'
$ws.Range('D2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.299.32'
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.03%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '514.06'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +4.06%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '130.71'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +2.01%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.995'
$ws.Range('D7').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.533'
$ws.Range('D8').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.81%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.319.01'
$ws.Range('D9').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.97%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.103'
$ws.Range('D10').Style = "Normal"

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +8.85%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.10'
$ws.Range('D12').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +7.86%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.00'
$ws.Range('D14').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +5.98%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.712.64'
$ws.Range('D15').Style = "Normal"

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '55.912.45'
$ws.Range('D16').Style = "Normal"

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +3.78%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +4.21%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.289.46'
$ws.Range('D18').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.54'
$ws.Range('D19').Style = "Normal"

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +3.33%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '319.25'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +6.37%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +4.95%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.39'
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.994'
$ws.Range('D25').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +6.34%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.64'
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +5.24%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '170.74'
$ws.Range('D28').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +9.56%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.24'
$ws.Range('D30').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +5.46%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0720'
$ws.Range('D31').Style = "Normal"

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +5.39%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +2.59%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.991'
$ws.Range('D35').Style = "Normal"

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +5.37%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +7.27%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'Stacks'
$ws.Range('B39').Style = "Normal"

$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C39').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.51'
$ws.Range('D39').Style = "Normal"

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +8.99%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'OKB'
$ws.Range('B40').Style = "Normal"

$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C40').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '37.21'
$ws.Range('D40').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +3.90%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.59%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '138.06'
$ws.Range('D42').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +10.12%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.56'
$ws.Range('D43').Style = "Normal"

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +6.31%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.12'
$ws.Range('D44').Style = "Normal"

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +7.32%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '263.99'
$ws.Range('D45').Style = "Normal"

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +10.55%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +4.30%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +4.16%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.554'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'Polygon'
$ws.Range('B49').Style = "Normal"

$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('C49').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.382'
$ws.Range('D49').Style = "Normal"

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.64%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'VeChain'
$ws.Range('B50').Style = "Normal"

$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C50').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0215'
$ws.Range('D50').Style = "Normal"

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +5.72%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '16.77'
$ws.Range('D51').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +4.67%  '
$ws.Range('E51').Style = "Normal"

